$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously only had a combined label in B2 ("SAN DIEGO AREA TOTALS").
# Split it: A2 gets the area label, B2 becomes "Totals" (matching the pattern
# used for the other port sections further down the sheet, e.g. row 20).
$ws.Range("A2").Value = "SAN DIEGO AREA TOTALS"
$ws.Range("B2").Value = "Totals"

# Column A now holds the same kind of (wider) text as column B, so widen it
# to match column B's width.
$ws.Range("A1").EntireColumn.ColumnWidth = 22.5

# Update the active cell / selection shown when the sheet is opened.
$ws.Range("B3").Select() | Out-Null
